# Generate Report for Handback
# Update the recorded timestamps for the 7f2d55f9-1250-4f40-81b2-ae1e03b433b6
# handback entry across the Overview, zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 7f2d55f9-...md
$wsOverview.Range("G4").Value = "2017-02-17 08:13:58"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2017-02-17 08:13:42"
$wsZhCn.Range("L4").Value = "2017-02-17 08:14:43"

# de-de sheet: Correspond Handback DateTime
# (note: de-de's "Correspond Handoff Datetime" H4 shares the same value as
#  Overview G4, so updating Overview G4 above also updates this cell's
#  displayed value since they reference the same underlying shared string)
$wsDeDe.Range("L4").Value = "2017-02-17 08:15:06"
